$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "CU-usuarios" -> "CU-" | "UC1-" | "usuarios" (three runs) by
#    inserting "UC1-" right after "CU-". A zero-width InsertAfter
#    causes this engine to coalesce every adjacent run that shares
#    identical formatting (not just the two runs touched), so the
#    following "-iniciar " / "sesión" runs get folded in too. Toggling
#    Bold on/off (net no-op) on each logical segment forces the
#    engine to keep/restore separate <w:r> boundaries there.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("CU-usuarios") | Out-Null
$matchStart = $r.Start
$matchEnd = $r.End

$insertStart = $matchStart + 3  # right after "CU-"
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertAfter("UC1-")

$runCU = $d.Range($matchStart, $insertStart)
$runCU.Font.Bold = $true
$runCU.Font.Bold = $false

$runUC1 = $d.Range($insertStart, $insertStart + 4)
$runUC1.Font.Bold = $true
$runUC1.Font.Bold = $false

$runUsuarios = $d.Range($insertStart + 4, $matchEnd + 4)
$runUsuarios.Font.Bold = $true
$runUsuarios.Font.Bold = $false

$runIniciar = $d.Range($matchEnd + 4, $matchEnd + 4 + 9)
$runIniciar.Font.Bold = $true
$runIniciar.Font.Bold = $false

$runSesion = $d.Range($matchEnd + 4 + 9, $matchEnd + 4 + 9 + 6)
$runSesion.Font.Bold = $true
$runSesion.Font.Bold = $false

# ------------------------------------------------------------------
# 2) Merge the three runs describing the incorrect-password attempt
#    into a single run with the concatenated text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "El usuario ingresa el usuario del sistema y la contraseña incorrecta y da clic en ingresar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario ingresa el usuario del sistema y la contraseña incorrecta y da clic en ingresar",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) Merge the two runs describing the system's validation response.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "El sistema valida con la base de datos, avisa que la el usuario o la contraseña son incorrectos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema valida con la base de datos, avisa que la el usuario o la contraseña son incorrectos",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Merge the two runs describing the forgotten-password click.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "El usuario no recuerda la contraseña, da clic en ¿olvido su contraseña?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario no recuerda la contraseña, da clic en ¿olvido su contraseña?",
    2) | Out-Null
